$wb = $excel.ActiveWorkbook

# --- Text / shared-string content updates -------------------------------

# Sheet "Cont adminstrator": row 15 ("Educație" business)
$wsCont = $wb.Worksheets.Item("Cont adminstrator")
$wsCont.Range("B15").Value = "Scoala particulara Suryia"
$wsCont.Range("C15").Value = "suryoiascoalaa1@automation.33mail.com"

# Sheet "Receptie": receptionist emails
$wsRec = $wb.Worksheets.Item("Receptie")
$wsRec.Range("B2").Value = "loise3@staffcalendis.33mail.com"
$wsRec.Range("B3").Value = "kamceatka@staffcalendis.33mail.com"
$wsRec.Range("B4").Value = "groblins1@staffcalendis.33mail.com"

# Sheet "Angajati": employee emails
$wsAng = $wb.Worksheets.Item("Angajati")
$wsAng.Range("B2").Value = "republicak@staffcalendis.33mail.com"
$wsAng.Range("B3").Value = "fermancow@staffcalendis.33mail.com"
$wsAng.Range("B4").Value = "comabatfuoo@staffcalendis.33mail.com"
$wsAng.Range("B5").Value = "mobentulju@staffcalendis.33mail.com"

# --- Selection (active cell) updates ------------------------------------

$wsCont.Range("C15").Select()
$wsRec.Range("B4").Select()

# Restore the originally active sheet/tab ("Angajati") so the workbook's
# active tab is unaffected by the selection changes above.
$wsAng.Activate()

# --- New generation of hidden AutoFilter defined names -------------------
# Each autofiltered sheet accrues one more "_xlnm._FilterDatabase_0_..._0"
# defined name (local to the sheet). Reproduce the new 91st-generation name
# for each of the three autofiltered sheets.

$suffix = ""
for ($i = 0; $i -lt 90; $i++) {
    $suffix = $suffix + "_0"
}
$filterName = "_xlnm._FilterDatabase" + $suffix

$wsContFd = $wb.Worksheets.Item("Cont adminstrator")
$wsContFd.Names.Add($filterName, "='Cont adminstrator'!`$A`$1:`$A`$19")

$wsDomenii = $wb.Worksheets.Item("Domenii")
$wsDomenii.Names.Add($filterName, "=Domenii!`$A`$4:`$A`$7")

$wsDomeniiExist = $wb.Worksheets.Item("Domenii existente")
$wsDomeniiExist.Names.Add($filterName, "='Domenii existente'!`$A`$1:`$Q`$15")
